# DALA-4606: adapt quality buckets
# - quality bucket "Incomplete" has been replaced by "NoDataFound"
# - "NoEvidenceFound" has been removed as a value for Yes-No questions
#   (i.e. the dropdown option string "Yes/No/No Evidence Found" becomes "Yes/No")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All cells in column "Options" (G) whose value was the quality-bucket
# string "Yes/No/No Evidence Found" get migrated to the new "Yes/No" value.
$cellsToUpdate = @(
    "G21","G30","G31","G32","G33","G36","G37","G38","G39","G40",
    "G41","G42","G43","G44","G45","G46","G47","G48","G49","G50",
    "G51","G52","G59","G63","G74","G75","G76","G77","G78","G79",
    "G81","G82","G83","G84","G85","G88","G89","G90","G95","G96",
    "G97","G98","G99"
)

foreach ($cellRef in $cellsToUpdate) {
    $ws.Range($cellRef).Value = "Yes/No"
}

# Keep the current selection on the active sheet in line with where the
# edits were made.
$ws.Range("G40").Select()
